$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.369.22'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.516.44'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.31'
$ws.Range('E5').Value = '  +1.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.68'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.64'
$ws.Range('E9').Value = '  +7.10%  '
$ws.Range('E10').Value = '  +0.22%  '
$ws.Range('E11').Value = '  +4.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.114.70'
$ws.Range('E12').Value = '  +0.38%  '
$ws.Range('E13').Value = '  +1.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000182'
$ws.Range('E14').Value = '  +1.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.514.36'
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.83'
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.362.89'
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('E18').Value = '  +2.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.77'
$ws.Range('E19').Value = '  +3.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.57'
$ws.Range('E20').Value = '  -2.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '394.40'
$ws.Range('E21').Value = '  +2.85%  '
$ws.Range('E22').Value = '  +1.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.656.98'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.64'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000118'
$ws.Range('E27').Value = '  +2.79%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  -1.84%  '
$ws.Range('E30').Value = '  +1.58%  '
$ws.Range('E31').Value = '  +0.29%  '
$ws.Range('E32').Value = '  -6.71%  '
$ws.Range('E33').Value = '  +7.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.547.66'
$ws.Range('E34').Value = '  +0.69%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.43'
$ws.Range('E36').Value = '  -0.58%  '
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.97'
$ws.Range('E38').Value = '  +1.64%  '
$ws.Range('E39').Value = '  +1.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '166.96'
$ws.Range('E40').Value = '  +1.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0791'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.812'
$ws.Range('E42').Value = '  +0.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '25.74'
$ws.Range('E43').Value = '  -1.92%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.46'
$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.67'
$ws.Range('E46').Value = '  +3.19%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.18'
$ws.Range('E47').Value = '  -1.76%  '
$ws.Range('E48').Value = '  +0.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.393.98'
$ws.Range('E49').Value = '  -3.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.900'
$ws.Range('E50').Value = '  -2.08%  '
$ws.Range('E51').Value = '  +0.39%  '
